$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update rows 11-12: replace the "4. Identification as" / "5. Requirement for
# workplace accommodations" descriptions with the new Orientation / Accomodations rows.
$ws.Range("C11").Value = "Orientation: 'string',"
$ws.Range("C12").Value = "Accomodations: 'string',"

# Re-apply the quote-prefix text style (matching the other C3:C47 description cells)
# that gets cleared when the cell value is reassigned.
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C11:C12").PasteSpecial(-4122) | Out-Null

# --- Append new rows 48-52 for the URM / degree "no answer" / clinical categories.
$ws.Range("B48").Value = -5
$ws.Range("C48").Value = "urmNoAnswer"
$ws.Range("B49").Value = -4
$ws.Range("C49").Value = "degClinical"
$ws.Range("B50").Value = -3
$ws.Range("C50").Value = "degNonClinical"
$ws.Range("B51").Value = -2
$ws.Range("C51").Value = "degOther"
$ws.Range("B52").Value = -1
$ws.Range("C52").Value = "degNoAnswer"

# Match the numeric index column's style (centered/top-aligned) used by B2:B47.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B48:B52").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Match the saved selection state recorded in the workbook.
$ws.Range("C53").Select() | Out-Null
